# This script reproduces the "Updated cryptos list" data refresh:
# new Price (D) / Volume 1h change (E) figures, plus the NEARProtocol /
# PancakeSwap rows trading rank positions (31 <-> 32).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.757.94"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "3.440.15"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.18"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.67"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "3.438.82"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.91"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D13").Value = "4.034.78"
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.64"
$ws.Range("E15").Value = "  -7.33%  "
$ws.Range("D16").Value = "65.795.17"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "3.437.40"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.74"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "366.86"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.56"
$ws.Range("E22").Value = "  -2.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.00"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000122"
$ws.Range("E26").Value = "  +2.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.97"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.71"
$ws.Range("E32").Value = "  -2.82%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -4.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.95"
$ws.Range("E35").Value = "  -2.30%  "
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.66"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.878"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.52"
$ws.Range("E39").Value = "  +4.78%  "
$ws.Range("E40").Value = "  -1.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.60"
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("D42").Value = "2.748.53"
$ws.Range("E42").Value = "  +2.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.45"
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.43"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  -2.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.28"
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.64"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "326.18"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("E51").Value = "  +0.65%  "
